$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking values so they stay as text (matching inlineStr source)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.528.75"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.491.04"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "314.15"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "93.34"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "32.69"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "2.876.11"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "16.19"
$ws.Range("E15").Value = "  +10.33%  "
$ws.Range("D16").Value = "2.479.08"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "0.759"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "41.555.58"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "71.10"
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("D22").Value = "11.25"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "236.46"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "2.72"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "25.25"
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "36.19"
$ws.Range("D31").Value = "157.75"
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("D32").Value = "5.47"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("D33").Value = "2.58"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "17.80"
$ws.Range("E35").Value = "  +5.45%  "
$ws.Range("E36").Value = "  -6.48%  "
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "19.83"
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("D44").Value = "1.960.09"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").Value = "8.91"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").Value = "2.728.77"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").Value = "96.83"
$ws.Range("D50").Value = "67.88"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "73.92"

# Reset style so no extra formatting attribute lingers on the cell
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
